{"js": "// Week 15 Update 6.3\n// Rewrites the REST-best-practices answer and the \"other REST features\" answer,\n// and replaces the trailing \"<div></div>\" placeholder paragraph with a real answer.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the paragraphs we need by their current (pre-edit) text so the\n// script is resilient to exact indices.\nlet otherFeaturesPara = null;   // \"What are some other features the REST entails...\"\nlet favoriteThingPara = null;   // \"What is your favorite thing you learned this week?\"\nlet blankAfterFavoritePara = null; // blank paragraph right after the \"favorite thing\" question\nlet divPara = null;             // \"<div></div>\"\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text.trim();\n  if (t === \"What are some other features the REST entails that we didn't cover this week?\") {\n    otherFeaturesPara = items[i];\n  } else if (t === \"What is your favorite thing you learned this week?\") {\n    favoriteThingPara = items[i];\n    if (i + 1 < items.length && items[i + 1].text.trim() === \"\") {\n      blankAfterFavoritePara = items[i + 1];\n    }\n  } else if (t === \"<div></div>\") {\n    divPara = items[i];\n  }\n}\n\nif (!otherFeaturesPara || !favoriteThingPara || !divPara || !blankAfterFavoritePara) {\n  throw new Error(\"Could not locate expected paragraphs in the document.\");\n}\n\n// 1) Turn the old \"other features\" question into the first bullet-style answer,\n//    then insert the next two answers after it, then a blank paragraph, then\n//    re-insert the original \"other features\" question text as its own paragraph.\notherFeaturesPara.insertText(\n  \"Collection endpoints should use plural nouns, this helps remind and inform users the endpoint could/should have multiple entries.\",\n  Word.InsertLocation.replace\n);\n\nlet cursor = otherFeaturesPara.insertParagraph(\n  \"Use standard HTTP error codes for handling response/requests. This makes it easier to troubleshoot issues and narrow down resolutions. \",\n  Word.InsertLocation.after\n);\n\ncursor = cursor.insertParagraph(\n  \"Be clear with the naming convention on endpoints when it comes to versions, this way an individual can identify the version they need for the endpoint they need to access.\",\n  Word.InsertLocation.after\n);\n\ncursor = cursor.insertParagraph(\"\", Word.InsertLocation.after);\n\ncursor = cursor.insertParagraph(\n  \"What are some other features the REST entails that we didn't cover this week?\",\n  Word.InsertLocation.after\n);\n\n// 2) Remove the \"favorite thing\" question paragraph and the following blank\n//    paragraph (the \"<div></div>\" placeholder paragraph is kept and reused).\nfavoriteThingPara.delete();\nblankAfterFavoritePara.delete();\n\n// 3) Replace the \"<div></div>\" placeholder text with the real answer.\ndivPara.insertText(\n  \"There are multiple methods which were not covered, HEAD, PATCH, OPTIONS, and TRACE. There area also HEADER, NAVIGATION, EVENTING, Error codes, even Authentications features.\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Week 15 Update 6.3\n# Rewrites the REST-best-practices answer and the \"other REST features\" answer,\n# and replaces the trailing \"<div></div>\" placeholder paragraph with a real answer.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraphs we need by their current (pre-edit) text so the\n# script is resilient to exact indices.\n$otherFeaturesIndex = -1\n$favoriteThingIndex = -1\n$blankAfterFavoriteIndex = -1\n$divIndex = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text.Trim()\n    if ($t -eq \"What are some other features the REST entails that we didn't cover this week?\") {\n        $otherFeaturesIndex = $i\n    } elseif ($t -eq \"What is your favorite thing you learned this week?\") {\n        $favoriteThingIndex = $i\n        if (($i + 1) -le $d.Paragraphs.Count) {\n            $nextT = $d.Paragraphs($i + 1).Range.Text.Trim()\n            if ($nextT -eq \"\") {\n                $blankAfterFavoriteIndex = $i + 1\n            }\n        }\n    } elseif ($t -eq \"<div></div>\") {\n        $divIndex = $i\n    }\n}\n\nif ($otherFeaturesIndex -eq -1 -or $favoriteThingIndex -eq -1 -or $divIndex -eq -1 -or $blankAfterFavoriteIndex -eq -1) {\n    throw \"Could not locate expected paragraphs in the document.\"\n}\n\n# 1) Turn the old \"other features\" question into the first new answer, then\n#    insert the next two answers after it, then a blank paragraph, then\n#    re-insert the original \"other features\" question text as its own paragraph.\n$otherFeaturesPara = $d.Paragraphs($otherFeaturesIndex)\n$otherFeaturesPara.Range.Text = \"Collection endpoints should use plural nouns, this helps remind and inform users the endpoint could/should have multiple entries.\"\n\n$otherFeaturesPara.Range.InsertParagraphAfter()\n$newPara1 = $d.Paragraphs($otherFeaturesIndex + 1)\n$newPara1.Range.Text = \"Use standard HTTP error codes for handling response/requests. This makes it easier to troubleshoot issues and narrow down resolutions. \"\n\n$newPara1.Range.InsertParagraphAfter()\n$newPara2 = $d.Paragraphs($otherFeaturesIndex + 2)\n$newPara2.Range.Text = \"Be clear with the naming convention on endpoints when it comes to versions, this way an individual can identify the version they need for the endpoint they need to access.\"\n\n$newPara2.Range.InsertParagraphAfter()\n# leave the next paragraph (blank) untouched\n\n$blankPara = $d.Paragraphs($otherFeaturesIndex + 3)\n$blankPara.Range.InsertParagraphAfter()\n$newPara3 = $d.Paragraphs($otherFeaturesIndex + 4)\n$newPara3.Range.Text = \"What are some other features the REST entails that we didn't cover this week?\"\n\n# 2) Remove the \"favorite thing\" question paragraph and the following blank\n#    paragraph (the \"<div></div>\" placeholder paragraph is kept and reused).\n# Re-resolve indices: 4 new paragraphs were inserted after $otherFeaturesIndex,\n# so everything from the original $favoriteThingIndex onward shifted by +4.\n$shift = 4\n$favoriteThingPara = $d.Paragraphs($favoriteThingIndex + $shift)\n$blankAfterFavoritePara = $d.Paragraphs($blankAfterFavoriteIndex + $shift)\n$divPara = $d.Paragraphs($divIndex + $shift)\n\n$favoriteThingPara.Range.Delete()\n$blankAfterFavoritePara.Range.Delete()\n\n# 3) Replace the \"<div></div>\" placeholder text with the real answer.\n$divPara.Range.Text = \"There are multiple methods which were not covered, HEAD, PATCH, OPTIONS, and TRACE. There area also HEADER, NAVIGATION, EVENTING, Error codes, even Authentications features.\"\n"}
